$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32 (Leve Item ID 5484)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1425
$ws.Range("J32").Value = 1425
$ws.Range("L32").Value = 1425
$ws.Range("N32").Value = -2077

# Sheet ALC, row 51 (Leve Item ID 5486)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11000
$ws.Range("I51").Value = 12000
$ws.Range("J51").Value = 10666.667
$ws.Range("K51").Value = 12000
$ws.Range("L51").Value = 10666.667
$ws.Range("M51").Value = -11516
$ws.Range("N51").Value = -11634.667

# Sheet ALC, row 112 (Leve Item ID 27960)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1262.125
$ws.Range("J112").Value = 1428.1428
$ws.Range("L112").Value = 4284.428400000001
$ws.Range("N112").Value = -6500.428400000001

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5508.3335
$ws.Range("I132").Value = 2262.5
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 6787.5
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -4257.5
$ws.Range("N132").Value = -41060

# Sheet BSM, row 16 (Leve Item ID 1684)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1880
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""

# Sheet BSM, row 99 (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1225
$ws.Range("I99").Value = 1225
$ws.Range("K99").Value = 1225
$ws.Range("M99").Value = 273

# Sheet BSM, row 133 (Leve Item ID 43209)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120

# Sheet CRP, row 12 (Leve Item ID 1604)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 8599.799999999999
$ws.Range("J12").Value = 19999.5
$ws.Range("L12").Value = 19999.5
$ws.Range("N12").Value = -20339.5

# Sheet CRP, row 13 (Leve Item ID 1996)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 5650
$ws.Range("I13").Value = 5300
$ws.Range("J13").Value = 6000
$ws.Range("K13").Value = 5300
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = -5161
$ws.Range("N13").Value = -6278

# Sheet CRP, row 19 (Leve Item ID 2233)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 381
$ws.Range("I19").Value = 96.333336
$ws.Range("J19").Value = 665.6667
$ws.Range("K19").Value = 96.333336
$ws.Range("L19").Value = 665.6667
$ws.Range("M19").Value = 73.666664
$ws.Range("N19").Value = -1005.6667

# Sheet CRP, row 24 (Leve Item ID 2233)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 381
$ws.Range("I24").Value = 96.333336
$ws.Range("J24").Value = 665.6667
$ws.Range("K24").Value = 96.333336
$ws.Range("L24").Value = 665.6667
$ws.Range("M24").Value = 73.666664
$ws.Range("N24").Value = -1005.6667

# Sheet CRP, row 26 (Leve Item ID 2004)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 887.5
$ws.Range("J26").Value = 887.5
$ws.Range("L26").Value = 887.5
$ws.Range("N26").Value = -1461.5

# Sheet CRP, row 41 (Leve Item ID 1917)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = ""

# Sheet CUL, row 3 (Leve Item ID 44094)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9499.25
$ws.Range("I3").Value = 4332.3335
$ws.Range("J3").Value = 25000
$ws.Range("K3").Value = 12997.0005
$ws.Range("L3").Value = 75000
$ws.Range("M3").Value = -12885.0005
$ws.Range("N3").Value = -75224

# Sheet CUL, row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = -3224

# Sheet CUL, row 26 (Leve Item ID 4746)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 40
$ws.Range("I26").Value = 29.833334
$ws.Range("J26").Value = 101
$ws.Range("K26").Value = 89.50000199999999
$ws.Range("L26").Value = 303
$ws.Range("M26").Value = 198.499998
$ws.Range("N26").Value = -879

# Sheet CUL, row 38 (Leve Item ID 4860)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 2604.6365
$ws.Range("I38").Value = 3217.8
$ws.Range("K38").Value = 9653.400000000001
$ws.Range("M38").Value = -9306.400000000001

# Sheet CUL, row 39 (Leve Item ID 4712)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1266.6666

# Sheet CUL, row 55 (Leve Item ID 4733)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1213
$ws.Range("I55").Value = 700
$ws.Range("J55").Value = 1555
$ws.Range("K55").Value = 2100
$ws.Range("L55").Value = 4665
$ws.Range("M55").Value = -1923
$ws.Range("N55").Value = -5019

# Sheet CUL, row 96 (Leve Item ID 19816)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""

# Sheet CUL, row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2131
$ws.Range("J131").Value = 2103.7036
$ws.Range("L131").Value = 6311.110799999999
$ws.Range("N131").Value = -16391.1108

# Sheet CUL, row 135 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = ""
$ws.Range("N135").Value = -14070

# Sheet CUL, row 140 (Leve Item ID 44097)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1011.75
$ws.Range("I140").Value = 1011.75
$ws.Range("K140").Value = 3035.25
$ws.Range("M140").Value = 2144.75

# Sheet GSM, row 7 (Leve Item ID 4197)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 10000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -9888
$ws.Range("N7").Value = ""

# Sheet GSM, row 8 (Leve Item ID 4197)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 10000
$ws.Range("I8").Value = 10000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -9861
$ws.Range("N8").Value = ""

# Sheet GSM, row 13 (Leve Item ID 2443)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 569.2222
$ws.Range("I13").Value = 149.75
$ws.Range("J13").Value = 904.8
$ws.Range("K13").Value = 149.75
$ws.Range("L13").Value = 904.8
$ws.Range("M13").Value = -10.75
$ws.Range("N13").Value = -1182.8

# Sheet LTW, row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 566.3333
$ws.Range("I55").Value = 550
$ws.Range("J55").Value = 599
$ws.Range("K55").Value = 550
$ws.Range("L55").Value = 599
$ws.Range("M55").Value = -377
$ws.Range("N55").Value = -945

# Sheet LTW, row 100 (Leve Item ID 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1582
$ws.Range("I100").Value = 1582
$ws.Range("K100").Value = 1582
$ws.Range("M100").Value = -1041

# Sheet LTW, row 108 (Leve Item ID 25655)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 14950
$ws.Range("J108").Value = 14950
$ws.Range("L108").Value = 14950
$ws.Range("N108").Value = -22630

# Sheet LTW, row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""

# Sheet WVR, row 19 (Leve Item ID 2666)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = ""

# Sheet WVR, row 62 (Leve Item ID 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5467.3335
$ws.Range("I62").Value = 4701
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 4701
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -4077
$ws.Range("N62").Value = -8248

# Sheet WVR, row 65 (Leve Item ID 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5467.3335
$ws.Range("I65").Value = 4701
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 23505
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -20385
$ws.Range("N65").Value = -41240

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1496.6
$ws.Range("I132").Value = 1496.6
$ws.Range("K132").Value = 4489.799999999999
$ws.Range("M132").Value = -1959.799999999999

# Sheet WVR, row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450
